$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UK")

# Allow for social care provision (it_*.txt rows) and alignment with social
# care receipt: insert 4 new rows at the top of the UK policy schedule
# table (above the existing uk_2011_std.txt row) for the Italy-style
# filenames that now also drive the UK sheet, shifting the rest of the
# table down by 4 rows.
$ws.Rows("2:5").Insert()

$ws.Cells.Item(2, 1).Value = "it_2015_std.txt"
$ws.Cells.Item(2, 2).Value = ""
$ws.Cells.Item(2, 3).Value = ""
$ws.Cells.Item(2, 4).Value = ""

$ws.Cells.Item(3, 1).Value = "it_2016_std.txt"
$ws.Cells.Item(3, 2).Value = ""
$ws.Cells.Item(3, 3).Value = ""
$ws.Cells.Item(3, 4).Value = ""

$ws.Cells.Item(4, 1).Value = "it_2020_web_std.txt"
$ws.Cells.Item(4, 2).Value = ""
$ws.Cells.Item(4, 3).Value = ""
$ws.Cells.Item(4, 4).Value = ""

$ws.Cells.Item(5, 1).Value = "it_2020_web_std_v2.txt"
$ws.Cells.Item(5, 2).Value = ""
$ws.Cells.Item(5, 3).Value = ""
$ws.Cells.Item(5, 4).Value = ""
